# Apply corrected financial data for 미원상사 (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: corrected figures
$ws.Range("D2").Value = 2472
$ws.Range("E2").Value = 177
$ws.Range("F2").Value = 177
$ws.Range("G2").Value = 172
$ws.Range("H2").Value = 117
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 18
$ws.Range("K2").Value = 2286
$ws.Range("L2").Value = 834
$ws.Range("M2").Value = 1453
$ws.Range("N2").Value = 1297
$ws.Range("O2").Value = 155
$ws.Range("P2").Value = 78
$ws.Range("Q2").Value = 179
$ws.Range("R2").Value = -165
$ws.Range("S2").Value = -12
$ws.Range("T2").Value = 156
$ws.Range("U2").Value = 24
$ws.Range("V2").Value = 468
$ws.Range("W2").Value = 7.17
$ws.Range("X2").Value = 4.72
$ws.Range("Y2").Value = 7.81
$ws.Range("Z2").Value = 5.13
$ws.Range("AA2").Value = 57.38
$ws.Range("AB2").Value = 1504.32
$ws.Range("AC2").Value = 1692
$ws.Range("AD2").Value = 14.3
$ws.Range("AE2").Value = 22212
$ws.Range("AF2").Value = 1.09
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 5864477

# Row 3: corrected figures
$ws.Range("D3").Value = 2394
$ws.Range("E3").Value = 185
$ws.Range("F3").Value = 185
$ws.Range("G3").Value = 214
$ws.Range("H3").Value = 170
$ws.Range("I3").Value = 143
$ws.Range("J3").Value = 27
$ws.Range("K3").Value = 2312
$ws.Range("L3").Value = 676
$ws.Range("M3").Value = 1636
$ws.Range("N3").Value = 1461
$ws.Range("O3").Value = 175
$ws.Range("P3").Value = 79
$ws.Range("Q3").Value = 364
$ws.Range("R3").Value = -117
$ws.Range("S3").Value = -216
$ws.Range("T3").Value = 115
$ws.Range("U3").Value = 250
$ws.Range("V3").Value = 266
$ws.Range("W3").Value = 7.71
$ws.Range("X3").Value = 7.09
$ws.Range("Y3").Value = 10.36
$ws.Range("Z3").Value = 7.38
$ws.Range("AA3").Value = 41.34
$ws.Range("AB3").Value = 1675.21
$ws.Range("AC3").Value = 2455
$ws.Range("AD3").Value = 11.27
$ws.Range("AE3").Value = 25137
$ws.Range("AF3").Value = 1.1
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5831842

# Row 4: corrected figures
$ws.Range("D4").Value = 2430
$ws.Range("E4").Value = 167
$ws.Range("F4").Value = 228
$ws.Range("G4").Value = 260
$ws.Range("H4").Value = 264
$ws.Range("I4").Value = 219
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 3296
$ws.Range("L4").Value = 904
$ws.Range("M4").Value = 2392
$ws.Range("N4").Value = 1835
$ws.Range("O4").Value = 557
$ws.Range("P4").Value = 79
$ws.Range("Q4").Value = 305
$ws.Range("R4").Value = -133
$ws.Range("S4").Value = -64
$ws.Range("T4").Value = 201
$ws.Range("U4").Value = 103
$ws.Range("V4").Value = 318
$ws.Range("W4").Value = 6.87
$ws.Range("X4").Value = 10.85
$ws.Range("Y4").Value = 13.26
$ws.Range("Z4").Value = 9.4
$ws.Range("AA4").Value = 37.77
$ws.Range("AB4").Value = 2061.51
$ws.Range("AC4").Value = 3808
$ws.Range("AD4").Value = 8.48
$ws.Range("AE4").Value = 32161
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 5723939

# Row 5: corrected figures
$ws.Range("D5").Value = 3187
$ws.Range("E5").Value = 110
$ws.Range("F5").Value = 110
$ws.Range("G5").Value = 311
$ws.Range("H5").Value = 299
$ws.Range("I5").Value = 272
$ws.Range("J5").Value = 27
$ws.Range("K5").Value = 2901
$ws.Range("L5").Value = 829
$ws.Range("M5").Value = 2072
$ws.Range("N5").Value = 1811
$ws.Range("O5").Value = 261
$ws.Range("P5").Value = 80
$ws.Range("Q5").Value = 306
$ws.Range("R5").Value = -233
$ws.Range("S5").Value = -221
$ws.Range("T5").Value = 255
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 291
$ws.Range("W5").Value = 3.45
$ws.Range("X5").Value = 9.380000000000001
$ws.Range("Y5").Value = 14.93
$ws.Range("Z5").Value = 9.640000000000001
$ws.Range("AA5").Value = 40.03
$ws.Range("AB5").Value = 2132.69
$ws.Range("AC5").Value = 4878
$ws.Range("AD5").Value = 7.11
$ws.Range("AE5").Value = 35278
$ws.Range("AF5").Value = 0.98
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 5151757

# Row 6: corrected figures
$ws.Range("D6").Value = 3286
$ws.Range("E6").Value = 272
$ws.Range("F6").Value = 272
$ws.Range("G6").Value = 260
$ws.Range("H6").Value = 242
$ws.Range("I6").Value = 215
$ws.Range("K6").Value = 2455
$ws.Range("L6").Value = 479
$ws.Range("M6").Value = 1977
$ws.Range("N6").Value = 1977
$ws.Range("P6").Value = 80
$ws.Range("Q6").Value = 374
$ws.Range("R6").Value = -158
$ws.Range("S6").Value = -11
$ws.Range("T6").Value = 310
$ws.Range("U6").Value = 64
$ws.Range("V6").Value = 3
$ws.Range("W6").Value = 8.27
$ws.Range("X6").Value = 7.37
$ws.Range("Y6").Value = 11.37
$ws.Range("Z6").Value = 9.039999999999999
$ws.Range("AA6").Value = 24.21
$ws.Range("AB6").Value = 2326.05
$ws.Range("AC6").Value = 4184
$ws.Range("AD6").Value = 9.470000000000001
$ws.Range("AE6").Value = 39265
$ws.Range("AF6").Value = 1.01
$ws.Range("AG6").Value = 770
$ws.Range("AH6").Value = 1.94
$ws.Range("AI6").Value = 17.98
$ws.Range("AJ6").Value = 5034823

# Rows 7-9: remove stale/erroneous financial data, keep only index (A), period (B) and label (C)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
